# Apply cryptos.xlsx price/volume updates (commit: 'Updated symbol list on Sat Jan  7 23:31:56 UTC 2023 with GitHub Actions')
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

# Row 2
Set-TextValue $ws.Range("D2") "261.29"
Set-TextValue $ws.Range("E2") "0.71%"

# Row 3
Set-TextValue $ws.Range("D3") "27.15"
Set-TextValue $ws.Range("E3") "0.94%"

# Row 5
Set-TextValue $ws.Range("D5") "0.06205"
Set-TextValue $ws.Range("E5") "2.27%"

# Row 6
Set-TextValue $ws.Range("D6") "6.729"
Set-TextValue $ws.Range("E6") "0.43%"

# Row 7
Set-TextValue $ws.Range("D7") "0.8511"
Set-TextValue $ws.Range("E7") "-1.06%"

# Row 8
Set-TextValue $ws.Range("D8") "0.9071"
Set-TextValue $ws.Range("E8") "-1.65%"

# Row 9
Set-TextValue $ws.Range("D9") "0.1404"
Set-TextValue $ws.Range("E9") "0.32%"

# Row 10
Set-TextValue $ws.Range("D10") "0.04744"
Set-TextValue $ws.Range("E10") "-10.33%"

# Row 11
Set-TextValue $ws.Range("E11") "-0.22%"

# Row 12
Set-TextValue $ws.Range("D12") "0.03174"
Set-TextValue $ws.Range("E12") "1.30%"

# Row 13
Set-TextValue $ws.Range("D13") "0.09058"
Set-TextValue $ws.Range("E13") "-0.93%"

# Row 14
Set-TextValue $ws.Range("D14") "0.001536"
Set-TextValue $ws.Range("E14") "0.24%"

# Row 15
Set-TextValue $ws.Range("D15") "0.0006151"
Set-TextValue $ws.Range("E15") "1.73%"

# Row 16
Set-TextValue $ws.Range("D16") "0.006046"
Set-TextValue $ws.Range("E16") "-0.25%"

# Row 17
Set-TextValue $ws.Range("D17") "3.468"
Set-TextValue $ws.Range("E17") "-0.46%"

# Row 18
Set-TextValue $ws.Range("D18") "3.172"
Set-TextValue $ws.Range("E18") "0.04%"

# Row 19
Set-TextValue $ws.Range("E19") "-0.33%"

# Row 21
Set-TextValue $ws.Range("D21") "0.1291"
Set-TextValue $ws.Range("E21") "-0.52%"

# Row 22
Set-TextValue $ws.Range("D22") "4.122"
Set-TextValue $ws.Range("E22") "0.72%"

# Row 23
Set-TextValue $ws.Range("D23") "0.04221"
Set-TextValue $ws.Range("E23") "-0.27%"

# Row 24
Set-TextValue $ws.Range("E24") "0.19%"

# Row 25
Set-TextValue $ws.Range("D25") "0.004116"
Set-TextValue $ws.Range("E25") "1.84%"

# Row 26
Set-TextValue $ws.Range("E26") "0.14%"

# Row 40
Set-TextValue $ws.Range("E40") "0.75%"

# Row 41
Set-TextValue $ws.Range("D41") "0.1112"
Set-TextValue $ws.Range("E41") "-0.55%"

# Row 42
Set-TextValue $ws.Range("E42") "0.05%"

# Row 43
Set-TextValue $ws.Range("E43") "-0.68%"

# Row 44
Set-TextValue $ws.Range("D44") "0.01343"
Set-TextValue $ws.Range("E44") "-10.04%"

# Row 45
Set-TextValue $ws.Range("D45") "0.00005174"
Set-TextValue $ws.Range("E45") "-4.33%"

# Row 46
Set-TextValue $ws.Range("E46") "0.13%"

# Row 47
Set-TextValue $ws.Range("D47") "0.03591"
Set-TextValue $ws.Range("E47") "-34.14%"

# Row 48
Set-TextValue $ws.Range("D48") "0.05783"
Set-TextValue $ws.Range("E48") "-57.26%"

# Row 49
Set-TextValue $ws.Range("E49") "0.13%"

# Row 50
Set-TextValue $ws.Range("E50") "0.13%"
